# LV_Activities Changes - 26 June 2024
#
# - "Users" sheet: the test user "James Craven" is replaced with
#   "Amanda Donovan"; the sheet's selection moves to A4.
# - "Company" sheet: the module value "FR Capital Provider" is replaced
#   with "ActivityCompany"; the sheet's selection moves to A2, and this
#   sheet becomes the active tab (replacing "UpdateActivity").

$wb = $excel.ActiveWorkbook

# --- Users sheet -----------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A2").Value = "Amanda Donovan"
$wsUsers.Range("A4").Select()

# --- Company sheet -----------------------------------------------------
$wsCompany = $wb.Worksheets.Item("Company")
$wsCompany.Range("A2").Value = "ActivityCompany"
$wsCompany.Range("A2").Select()
